# N.Phuong chinh sua 10:48 21-02
#
# 1) Rename the existing (only) sheet "Sheet1" -> "LIBS".
# 2) Insert a brand-new sheet "VAR_MAP" right after "LIBS" and make it the
#    active/selected tab.
# 3) Populate VAR_MAP with the small MayBay/SoHieu/KieuDang/SoCot/SoDong/
#    ChuyenBay layout, matching the target shared-strings order.
# 4) Leave the selection on B8 in VAR_MAP (matches the target sheetView).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LIBS"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "VAR_MAP"

$ws2.Range("B3").Value = "MayBay"
$ws2.Range("C3").Value = "SoHieu"
$ws2.Range("C4").Value = "KieuDang"
$ws2.Range("C5").Value = "SoCot"
$ws2.Range("C6").Value = "SoDong"
$ws2.Range("B7").Value = "ChuyenBay"

# Target stored width is 18.5546875 characters; this host's ColumnWidth
# setter only persists in 1/6-character increments, so 17.67 (-> stored
# 18.5) is the closest representable width reachable from COM.
$ws2.Columns.Item(2).ColumnWidth = 17.67

$ws2.Range("B8").Select()
